$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "13÷5=2, 3"
$t.Cell(1, 2).Range.Text = "94÷2=47, 0"
$t.Cell(1, 3).Range.Text = "63÷3=21, 0"
$t.Cell(1, 4).Range.Text = "72÷8=9, 0"
$t.Cell(1, 5).Range.Text = "79÷3=26, 1"
$t.Cell(5, 1).Range.Text = "90÷3=30, 0"
$t.Cell(5, 2).Range.Text = "91÷7=13, 0"
$t.Cell(5, 3).Range.Text = "73÷7=10, 3"
$t.Cell(5, 4).Range.Text = "63÷4=15, 3"
$t.Cell(5, 5).Range.Text = "59÷8=7, 3"
$t.Cell(9, 1).Range.Text = "34÷8=4, 2"
$t.Cell(9, 2).Range.Text = "41÷2=20, 1"
$t.Cell(9, 3).Range.Text = "24÷6=4, 0"
$t.Cell(9, 4).Range.Text = "35÷5=7, 0"
$t.Cell(9, 5).Range.Text = "33÷6=5, 3"
$t.Cell(13, 1).Range.Text = "40÷4=10, 0"
$t.Cell(13, 2).Range.Text = "31÷4=7, 3"
$t.Cell(13, 3).Range.Text = "47÷5=9, 2"
$t.Cell(13, 4).Range.Text = "40÷4=10, 0"
$t.Cell(13, 5).Range.Text = "92÷9=10, 2"
$t.Cell(17, 1).Range.Text = "83÷4=20, 3"
$t.Cell(17, 2).Range.Text = "65÷5=13, 0"
$t.Cell(17, 3).Range.Text = "85÷9=9, 4"
$t.Cell(17, 4).Range.Text = "77÷9=8, 5"
$t.Cell(17, 5).Range.Text = "10÷3=3, 1"
